$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ E = 3; G = 111.005114; H = 333.015342; I = 0.07617677061978854; J = 0.07617677061978853; K = 3; M = 49.32494666666667; N = 147.97484; O = 0.5781391668179461; P = 0.5781391668179461; Q = 5475.321327777254; R = 49277.89194999528; S = 0.04404077469700635; T = 0.04404077469700633 }
    3 = @{ E = 3; G = 111.005114; H = 333.015342; I = 0.07617677061978854; J = 0.07617677061978853; K = 3; M = 1.693583; N = 5.080749; O = 0.01985053671064022; P = 0.01985053671064022; Q = 187.996373983462; R = 1691.967365851158; S = 0.001512149781686132; T = 0.001512149781686132 }
    4 = @{ E = 3; G = 111.005114; H = 333.015342; I = 0.07617677061978854; J = 0.07617677061978853; K = 3; M = 27.10902833333334; N = 81.32708500000001; O = 0.3177457273252148; P = 0.3177457273252148; Q = 3009.240780570897; R = 27083.16702513807; S = 0.02420484338587076; T = 0.02420484338587076 }
    5 = @{ E = 3; G = 111.005114; H = 333.015342; I = 0.07617677061978854; J = 0.07617677061978853; K = 3; M = 7.189177999999999; N = 21.567534; O = 0.08426456914619893; P = 0.08426456914619893; Q = 798.0355234562919; R = 7182.319711106628; S = 0.006419002755225307; T = 0.006419002755225305 }
    6 = @{ E = 3; G = 930.1503093333332; H = 2790.450928; I = 0.6383115534299679; J = 0.6383115534299678; K = 3; M = 49.32494666666667; N = 147.97484; O = 0.5781391668179461; P = 0.5781391668179461; Q = 45879.61439985016; R = 412916.5295986515; S = 0.3690329096702705; T = 0.3690329096702705 }
    7 = @{ E = 3; G = 930.1503093333332; H = 2790.450928; I = 0.6383115534299679; J = 0.6383115534299678; K = 3; M = 1.693583; N = 5.080749; O = 0.01985053671064022; P = 0.01985053671064022; Q = 1575.286751331675; R = 14177.58076198507; S = 0.01267082692418736; T = 0.01267082692418736 }
    8 = @{ E = 3; G = 930.1503093333332; H = 2790.450928; I = 0.6383115534299679; J = 0.6383115534299678; K = 3; M = 27.10902833333334; N = 81.32708500000001; O = 0.3177457273252148; P = 0.3177457273252148; Q = 25215.4710899761; R = 226939.2398097849; S = 0.2028207688046929; T = 0.2028207688046928 }
    9 = @{ E = 3; G = 930.1503093333332; H = 2790.450928; I = 0.6383115534299679; J = 0.6383115534299678; K = 3; M = 7.189177999999999; N = 21.567534; O = 0.08426456914619893; P = 0.08426456914619893; Q = 6687.016140552393; R = 60183.14526497154; S = 0.05378704803081718; T = 0.05378704803081717 }
    10 = @{ E = 3; G = 72.47344466666667; H = 217.420334; I = 0.04973458223193756; J = 0.04973458223193756; K = 3; M = 49.32494666666667; N = 147.97484; O = 0.5781391668179461; P = 0.5781391668179461; Q = 3574.748792932951; R = 32172.73913639656; S = 0.028753509933611; T = 0.028753509933611 }
    11 = @{ E = 3; G = 72.47344466666667; H = 217.420334; I = 0.04973458223193756; J = 0.04973458223193756; K = 3; M = 1.693583; N = 5.080749; O = 0.01985053671064022; P = 0.01985053671064022; Q = 122.7397938389073; R = 1104.658144550166; S = 0.0009872581503834312; T = 0.0009872581503834312 }
    12 = @{ E = 3; G = 72.47344466666667; H = 217.420334; I = 0.04973458223193756; J = 0.04973458223193756; K = 3; M = 27.10902833333334; N = 81.32708500000001; O = 0.3177457273252148; P = 0.3177457273252148; Q = 1964.684664882933; R = 17682.16198394639; S = 0.0158029510045027; T = 0.0158029510045027 }
    13 = @{ E = 3; G = 72.47344466666667; H = 217.420334; I = 0.04973458223193756; J = 0.04973458223193756; K = 3; M = 7.189177999999999; N = 21.567534; O = 0.08426456914619893; P = 0.08426456914619893; Q = 521.0244939818173; R = 4689.220445836356; S = 0.004190863143440419; T = 0.004190863143440419 }
    14 = @{ E = 3; G = 343.5753833333333; H = 1030.72615; I = 0.2357770937183061; J = 0.235777093718306; K = 3; M = 49.32494666666667; N = 147.97484; O = 0.5781391668179461; P = 0.5781391668179461; Q = 16946.83745889622; R = 152521.537130066; S = 0.1363119725170583; T = 0.1363119725170582 }
    15 = @{ E = 3; G = 343.5753833333333; H = 1030.72615; I = 0.2357770937183061; J = 0.235777093718306; K = 3; M = 1.693583; N = 5.080749; O = 0.01985053671064022; P = 0.01985053671064022; Q = 581.8734284318167; R = 5236.86085588635; S = 0.004680301854383294; T = 0.004680301854383294 }
    16 = @{ E = 3; G = 343.5753833333333; H = 1030.72615; I = 0.2357770937183061; J = 0.235777093718306; K = 3; M = 27.10902833333334; N = 81.32708500000001; O = 0.3177457273252148; P = 0.3177457273252148; Q = 9313.994801419196; R = 83825.95321277276; S = 0.0749171641301485; T = 0.0749171641301485 }
    17 = @{ E = 3; G = 343.5753833333333; H = 1030.72615; I = 0.2357770937183061; J = 0.235777093718306; K = 3; M = 7.189177999999999; N = 21.567534; O = 0.08426456914619893; P = 0.08426456914619893; Q = 2470.024587201566; R = 22230.2212848141; S = 0.01986765521671603; T = 0.01986765521671603 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}